$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the newly added columns (BG:BL) from column A of the same row
foreach ($r in 2..12) {
    $ws.Range("A" + $r).Copy() | Out-Null
    $ws.Range("BG" + $r + ":BL" + $r).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Header row (row 1): rename / add headers for columns AQ:BL
$ws.Range("AQ1").Value = "percent_women"
$ws.Range("AR1").Value = "cuij_responders_active"
$ws.Range("AS1").Value = "cuij_responders_control"
$ws.Range("AT1").Value = "comorbid_mental."
$ws.Range("AU1").Value = "country"
$ws.Range("AV1").Value = "age_m_active"
$ws.Range("AW1").Value = "age_sd_active"
$ws.Range("AX1").Value = "age_m_control"
$ws.Range("AY1").Value = "age_sd_control"
$ws.Range("AZ1").Value = "age_m_overall"
$ws.Range("BA1").Value = "age_sd_overall"
$ws.Range("BB1").Value = "active_percent_women"
$ws.Range("BC1").Value = "control_percent_women"
$ws.Range("BD1").Value = "Notes"
$ws.Range("BE1").Value = "comorbid_mental?"
$ws.Range("BF1").Value = "overall_percent_women"
$ws.Range("BG1").Value = "overall_mean_age"
$ws.Range("BH1").Value = "mean_age"
$ws.Range("BI1").Value = "pooled_sd_age"
$ws.Range("BJ1").Value = "sd_age"
$ws.Range("BK1").Value = "cohens_d_active"
$ws.Range("BL1").Value = "cohens_d_control"

# Year column (B) - force text values with trailing .0
$ws.Range("B6").Value = "'1990.0"
$ws.Range("B9").Value = "'2009.0"
$ws.Range("B10").Value = "'2006.0"
$ws.Range("B11").Value = "'2004.0"
$ws.Range("B12").Value = "'2006.0"

# Data cells in columns AQ:BL (rows 2-12)
$ws.Range("AS2").Value = ""
$ws.Range("AT2").Value = "missing"
$ws.Range("AU2").Value = "USA"
$ws.Range("BA2").Value = ""
$ws.Range("BD2").Value = "FDA Trial - Poor reporting"
$ws.Range("BG2").Value = ""
$ws.Range("BH2").Value = ""
$ws.Range("BI2").Value = ""
$ws.Range("BJ2").Value = ""
$ws.Range("BK2").Value = ""
$ws.Range("BL2").Value = ""
$ws.Range("AS3").Value = ""
$ws.Range("AT3").Value = "missing"
$ws.Range("AU3").Value = "USA"
$ws.Range("BG3").Value = ""
$ws.Range("BH3").Value = ""
$ws.Range("BI3").Value = ""
$ws.Range("BJ3").Value = ""
$ws.Range("BK3").Value = ""
$ws.Range("BL3").Value = ""
$ws.Range("AT4").Value = ""
$ws.Range("AU4").Value = "USA"
$ws.Range("BG4").Value = ""
$ws.Range("BH4").Value = ""
$ws.Range("BI4").Value = ""
$ws.Range("BJ4").Value = ""
$ws.Range("BK4").Value = ""
$ws.Range("BL4").Value = ""
$ws.Range("AQ5").Value = 45.5044510385757
$ws.Range("AT5").Value = ""
$ws.Range("AU5").Value = "USA"
$ws.Range("AW5").Value = 2.6
$ws.Range("AX5").Value = 12.3
$ws.Range("AY5").Value = 2.6
$ws.Range("AZ5").Value = ""
$ws.Range("BB5").Value = 44
$ws.Range("BC5").Value = 47
$ws.Range("BD5").Value = ""
$ws.Range("BF5").Value = 45.5044510385757
$ws.Range("BG5").Value = 12.2501483679525
$ws.Range("BH5").Value = 12.2501483679525
$ws.Range("BI5").Value = 2.6
$ws.Range("BJ5").Value = 2.6
$ws.Range("BK5").Value = ""
$ws.Range("BL5").Value = ""
$ws.Range("AS6").Value = ""
$ws.Range("AT6").Value = "Anxiety-related disorders (Group1) and Conduct or oppositional disorders (group2)"
$ws.Range("BG6").Value = ""
$ws.Range("BH6").Value = ""
$ws.Range("BI6").Value = ""
$ws.Range("BJ6").Value = ""
$ws.Range("BK6").Value = ""
$ws.Range("BL6").Value = ""
$ws.Range("AQ7").Value = 50.8126984126984
$ws.Range("AU7").Value = ""
$ws.Range("AW7").Value = ""
$ws.Range("AX7").Value = 12.4
$ws.Range("AY7").Value = ""
$ws.Range("AZ7").Value = ""
$ws.Range("BA7").Value = ""
$ws.Range("BB7").Value = 47.6
$ws.Range("BC7").Value = 56.8
$ws.Range("BD7").Value = "FDA trial - poor reporting"
$ws.Range("BF7").Value = 50.8126984126984
$ws.Range("BG7").Value = 12.3349206349206
$ws.Range("BH7").Value = 12.3349206349206
$ws.Range("BI7").Value = ""
$ws.Range("BJ7").Value = ""
$ws.Range("BK7").Value = ""
$ws.Range("BL7").Value = ""
$ws.Range("AQ8").Value = 50.1730769230769
$ws.Range("AU8").Value = ""
$ws.Range("AW8").Value = ""
$ws.Range("AX8").Value = 12.3
$ws.Range("AY8").Value = ""
$ws.Range("AZ8").Value = ""
$ws.Range("BB8").Value = 50
$ws.Range("BC8").Value = 50.5
$ws.Range("BD8").Value = ""
$ws.Range("BF8").Value = 50.1730769230769
$ws.Range("BG8").Value = 12.0384615384615
$ws.Range("BH8").Value = 12.0384615384615
$ws.Range("BI8").Value = ""
$ws.Range("BJ8").Value = ""
$ws.Range("BK8").Value = ""
$ws.Range("BL8").Value = ""
$ws.Range("AQ9").Value = 60.7446428571429
$ws.Range("AT9").Value = ""
$ws.Range("AU9").Value = "Japan"
$ws.Range("AW9").Value = 1.99
$ws.Range("AX9").Value = 14.8
$ws.Range("AY9").Value = 2.62
$ws.Range("AZ9").Value = ""
$ws.Range("BB9").Value = 55.2
$ws.Range("BC9").Value = 66.7
$ws.Range("BD9").Value = ""
$ws.Range("BF9").Value = 60.7446428571429
$ws.Range("BG9").Value = 14.5928571428571
$ws.Range("BH9").Value = 14.5928571428571
$ws.Range("BI9").Value = 2.31483620730856
$ws.Range("BJ9").Value = 2.31483620730856
$ws.Range("BK9").Value = ""
$ws.Range("BL9").Value = ""
$ws.Range("AT10").Value = ""
$ws.Range("AU10").Value = "European, multi-centered (31 recruitment sites)"
$ws.Range("AW10").Value = 1
$ws.Range("AX10").Value = 16
$ws.Range("AY10").Value = 1
$ws.Range("AZ10").Value = 16
$ws.Range("BA10").Value = 1
$ws.Range("BC10").Value = ""
$ws.Range("BD10").Value = "*Ask Charlotte to check paper"
$ws.Range("BG10").Value = 16
$ws.Range("BH10").Value = 16
$ws.Range("BI10").Value = 1
$ws.Range("BJ10").Value = 1
$ws.Range("BK10").Value = ""
$ws.Range("BL10").Value = ""
$ws.Range("AQ11").Value = 53.4350574712644
$ws.Range("AS11").Value = ""
$ws.Range("AT11").Value = "Dysthmia, enuresis"
$ws.Range("AU11").Value = "USA"
$ws.Range("AW11").Value = 3.1
$ws.Range("AX11").Value = 12.1
$ws.Range("AY11").Value = 2.8
$ws.Range("AZ11").Value = 12.1
$ws.Range("BB11").Value = 52.8
$ws.Range("BC11").Value = 54.1
$ws.Range("BD11").Value = ""
$ws.Range("BF11").Value = 53.4350574712644
$ws.Range("BG11").Value = 12.1
$ws.Range("BH11").Value = 12.1
$ws.Range("BI11").Value = 2.95729291673125
$ws.Range("BJ11").Value = 2.95729291673125
$ws.Range("BK11").Value = ""
$ws.Range("BL11").Value = ""
$ws.Range("AQ12").Value = 51.9
$ws.Range("AS12").Value = ""
$ws.Range("AT12").Value = "Anxiety disorders (e.g. GAD, panic disorder, social phobia, specific phobia)"
$ws.Range("AU12").Value = "USA"
$ws.Range("AV12").Value = 12.2
$ws.Range("AW12").Value = 2.9
$ws.Range("AX12").Value = 12.4
$ws.Range("AY12").Value = 3
$ws.Range("AZ12").Value = 12.3
$ws.Range("BA12").Value = 3
$ws.Range("BB12").Value = 51.9
$ws.Range("BC12").Value = 51.9
$ws.Range("BD12").Value = ""
$ws.Range("BF12").Value = 51.9
$ws.Range("BG12").Value = 12.3007575757576
$ws.Range("BH12").Value = 12.3
$ws.Range("BI12").Value = 2.95080529828671
$ws.Range("BJ12").Value = 3
$ws.Range("BK12").Value = ""
$ws.Range("BL12").Value = ""
